$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 7860.1
$ws.Range("J17").Value = 8070.1377
$ws.Range("L17").Value = 24210.4131
$ws.Range("N17").Value = -24546.4131

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3480.6875
$ws.Range("I32").Value = 2541.5
$ws.Range("K32").Value = 2541.5
$ws.Range("M32").Value = -2215.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 6272.4546
$ws.Range("J116").Value = 7982.909
$ws.Range("L116").Value = 7982.909
$ws.Range("N116").Value = -14866.909

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 40005090
$ws.Range("I132").Value = 43483530
$ws.Range("J132").Value = 2990
$ws.Range("K132").Value = 130450590
$ws.Range("L132").Value = 8970
$ws.Range("M132").Value = -130448060
$ws.Range("N132").Value = -14030

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9785.967000000001
$ws.Range("I32").Value = 5648
$ws.Range("J32").Value = 24746.309
$ws.Range("K32").Value = 5648
$ws.Range("L32").Value = 24746.309
$ws.Range("M32").Value = -5361
$ws.Range("N32").Value = -25320.309

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4643748
$ws.Range("I45").Value = 6540709.5
$ws.Range("K45").Value = 6540709.5
$ws.Range("M45").Value = -6540332.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("N55").Value = 0
$ws.Range("L55").ClearContents()
$ws.Range("M55").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 3336276.2
$ws.Range("I102").Value = 4169139.2
$ws.Range("J102").Value = 4824
$ws.Range("K102").Value = 4169139.2
$ws.Range("L102").Value = 4824
$ws.Range("M102").Value = -4167517.2
$ws.Range("N102").Value = -8068

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1208556.4
$ws.Range("I110").Value = 1389690.4
$ws.Range("J110").Value = 996.6667
$ws.Range("K110").Value = 1389690.4
$ws.Range("L110").Value = 996.6667
$ws.Range("M110").Value = -1387645.4
$ws.Range("N110").Value = -5086.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5558212.5
$ws.Range("I86").Value = 7145683.5
$ws.Range("K86").Value = 7145683.5
$ws.Range("M86").Value = -7144560.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 5558212.5
$ws.Range("I89").Value = 7145683.5
$ws.Range("K89").Value = 35728417.5
$ws.Range("M89").Value = -35722801.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 5902380.5
$ws.Range("I94").Value = 9093679
$ws.Range("J94").Value = 51666.5
$ws.Range("K94").Value = 9093679
$ws.Range("L94").Value = 51666.5
$ws.Range("M94").Value = -9093228
$ws.Range("N94").Value = -52568.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 28575426
$ws.Range("I99").Value = 71430570
$ws.Range("J99").Value = 5332.3335
$ws.Range("K99").Value = 71430570
$ws.Range("L99").Value = 5332.3335
$ws.Range("M99").Value = -71429072
$ws.Range("N99").Value = -8328.333500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 9500
$ws.Range("I6").Value = 9500
$ws.Range("K6").Value = 9500
$ws.Range("M6").Value = -9387

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 95.76922999999999
$ws.Range("I7").Value = 102.5
$ws.Range("J7").Value = 15
$ws.Range("K7").Value = 102.5
$ws.Range("L7").Value = 15
$ws.Range("M7").Value = 10.5
$ws.Range("N7").Value = -241

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 50
$ws.Range("I25").Value = 50
$ws.Range("K25").Value = 50
$ws.Range("M25").Value = 124

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 21910.365
$ws.Range("I31").Value = 2986.3076
$ws.Range("J31").Value = 28218.385
$ws.Range("K31").Value = 2986.3076
$ws.Range("L31").Value = 28218.385
$ws.Range("M31").Value = -2691.3076
$ws.Range("N31").Value = -28808.385

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 21910.365
$ws.Range("I34").Value = 2986.3076
$ws.Range("J34").Value = 28218.385
$ws.Range("K34").Value = 2986.3076
$ws.Range("L34").Value = 28218.385
$ws.Range("M34").Value = -2784.3076
$ws.Range("N34").Value = -28622.385

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 9723.4
$ws.Range("I41").Value = 9723.4
$ws.Range("K41").Value = 9723.4
$ws.Range("M41").Value = -9295.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 11500

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 17319.666
$ws.Range("J51").Value = 49959
$ws.Range("L51").Value = 49959
$ws.Range("N51").Value = -51431

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6364.1724
$ws.Range("I58").Value = 8243.125
$ws.Range("K58").Value = 8243.125
$ws.Range("M58").Value = -8040.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 32500
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 4609.5713
$ws.Range("I60").Value = 4609.5713
$ws.Range("K60").Value = 4609.5713
$ws.Range("M60").Value = -4098.5713

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 17319.666
$ws.Range("J61").Value = 49959
$ws.Range("L61").Value = 49959
$ws.Range("N61").Value = -50655

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2798.3333
$ws.Range("I62").Value = 2798.3333
$ws.Range("K62").Value = 2798.3333
$ws.Range("M62").Value = -2174.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 2798.3333
$ws.Range("I65").Value = 2798.3333
$ws.Range("K65").Value = 13991.6665
$ws.Range("M65").Value = -10871.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4046.6667
$ws.Range("I99").Value = 3361.375
$ws.Range("J99").Value = 5417.25
$ws.Range("K99").Value = 3361.375
$ws.Range("L99").Value = 5417.25
$ws.Range("M99").Value = -1863.375
$ws.Range("N99").Value = -8413.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value = 8166.6665
$ws.Range("I103").Value = 1980.75
$ws.Range("K103").Value = 1980.75
$ws.Range("M103").Value = -808.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 4046.6667
$ws.Range("I126").Value = 3361.375
$ws.Range("J126").Value = 5417.25
$ws.Range("K126").Value = 10084.125
$ws.Range("L126").Value = 16251.75
$ws.Range("M126").Value = -7614.125
$ws.Range("N126").Value = -21191.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 6364.1724
$ws.Range("I136").Value = 8243.125
$ws.Range("K136").Value = 24729.375
$ws.Range("M136").Value = -22179.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 48633.285
$ws.Range("I5").Value = 706.25
$ws.Range("J5").Value = 201999.8
$ws.Range("K5").Value = 2118.75
$ws.Range("L5").Value = 605999.3999999999
$ws.Range("M5").Value = -2006.75
$ws.Range("N5").Value = -606223.3999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 996.3333
$ws.Range("J68").Value = 1450
$ws.Range("L68").Value = 4350
$ws.Range("N68").Value = -5972

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 996.3333
$ws.Range("J71").Value = 1450
$ws.Range("L71").Value = 13050
$ws.Range("N71").Value = -21162

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3908.2727
$ws.Range("J113").Value = 1985.8667
$ws.Range("L113").Value = 5957.6001
$ws.Range("N113").Value = -10297.6001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 48633.285
$ws.Range("I135").Value = 706.25
$ws.Range("J135").Value = 201999.8
$ws.Range("K135").Value = 6356.25
$ws.Range("L135").Value = 1817998.2
$ws.Range("M135").Value = -3821.25
$ws.Range("N135").Value = -1823068.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1505.6666
$ws.Range("I140").Value = 1381.375
$ws.Range("K140").Value = 4144.125
$ws.Range("M140").Value = 1035.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 7578900
$ws.Range("I113").Value = 16667849
$ws.Range("K113").Value = 16667849
$ws.Range("M113").Value = -16665679

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 42117.816
$ws.Range("I22").Value = 69279.08
$ws.Range("K22").Value = 69279.08
$ws.Range("M22").Value = -68984.08

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 42117.816
$ws.Range("I27").Value = 69279.08
$ws.Range("K27").Value = 69279.08
$ws.Range("M27").Value = -69172.08

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2460.8667
$ws.Range("J68").Value = 3041.5
$ws.Range("L68").Value = 3041.5
$ws.Range("N68").Value = -4539.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2460.8667
$ws.Range("J71").Value = 3041.5
$ws.Range("L71").Value = 15207.5
$ws.Range("N71").Value = -22695.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2888.2222
$ws.Range("I100").Value = 1998.3334
$ws.Range("K100").Value = 1998.3334
$ws.Range("M100").Value = -1457.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5437.231
$ws.Range("I122").Value = 3235.182
$ws.Range("J122").Value = 7052.067
$ws.Range("K122").Value = 9705.545999999998
$ws.Range("L122").Value = 21156.201
$ws.Range("M122").Value = -7255.545999999998
$ws.Range("N122").Value = -26056.201

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7179.829
$ws.Range("I132").Value = 7910.7188
$ws.Range("K132").Value = 23732.1564
$ws.Range("M132").Value = -21202.1564

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3475.7778
$ws.Range("I96").Value = 4113.6665
$ws.Range("J96").Value = 2200
$ws.Range("K96").Value = 4113.6665
$ws.Range("L96").Value = 2200
$ws.Range("M96").Value = -2740.6665
$ws.Range("N96").Value = -4946

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2507.9697
$ws.Range("I122").Value = 2366.6667
$ws.Range("K122").Value = 7100.000100000001
$ws.Range("M122").Value = -4650.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 26606810
$ws.Range("I132").Value = 37042556
$ws.Range("J132").Value = 991791.6
$ws.Range("K132").Value = 111127668
$ws.Range("L132").Value = 2975374.8
$ws.Range("M132").Value = -111125138
$ws.Range("N132").Value = -2980434.8
